# The deck ships two embedded themes:
#   ppt/theme/theme1.xml -> bound to the (only) Slide Master  ("Integral" / "Red Violet" palette)
#   ppt/theme/theme2.xml -> bound to the Notes Master         ("Office Theme" / "Office" palette)
#
# The authored change swaps the two themes' contents: the slides should now
# render with the stock "Office" colour palette while the notes master keeps
# the palette that used to belong to the slides ("Red Violet"). Font scheme
# and format scheme (fills/lines/effects) are identical between the two
# themes, so the swap is purely a colour-scheme (clrScheme) change.
#
# PowerPoint's automation surface exposes the legacy/standard 12-slot theme
# colour scheme through SlideMaster.ColorScheme (and the modern
# Slide.ThemeColorScheme, which maps to the same underlying scheme), where
# each slot's .RGB is independently settable. Re-map every slot from the
# "Integral" palette to the "Office" palette so the slide theme (theme1.xml)
# ends up holding the colours that theme2.xml used to hold.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$colorScheme = $master.ColorScheme

function Set-SchemeColor($scheme, $index, $r, $g, $b) {
    $scheme.Colors($index).RGB = $r + ($g * 256) + ($b * 65536)
}

# Target palette = the "Office" theme colours (previously theme2.xml),
# applied in clrScheme document order: dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink.
Set-SchemeColor $colorScheme 1  0x00 0x00 0x00   # dk1       000000
Set-SchemeColor $colorScheme 2  0xFF 0xFF 0xFF   # lt1       FFFFFF
Set-SchemeColor $colorScheme 3  0x44 0x54 0x6A   # dk2       44546A
Set-SchemeColor $colorScheme 4  0xE7 0xE6 0xE6   # lt2       E7E6E6
Set-SchemeColor $colorScheme 5  0x5B 0x9B 0xD5   # accent1   5B9BD5
Set-SchemeColor $colorScheme 6  0xED 0x7D 0x31   # accent2   ED7D31
Set-SchemeColor $colorScheme 7  0xA5 0xA5 0xA5   # accent3   A5A5A5
Set-SchemeColor $colorScheme 8  0xFF 0xC0 0x00   # accent4   FFC000
Set-SchemeColor $colorScheme 9  0x44 0x72 0xC4   # accent5   4472C4
Set-SchemeColor $colorScheme 10 0x70 0xAD 0x47   # accent6   70AD47
Set-SchemeColor $colorScheme 11 0x05 0x63 0xC1   # hlink     0563C1
Set-SchemeColor $colorScheme 12 0x95 0x4F 0x72   # folHlink  954F72
